$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Qui ?" (Who?) column (C) with the assigned person for each task
$ws.Range("C8").Value = "Alex"
$ws.Range("C9").Value = "Alex"
$ws.Range("C5").Value = "Fait"
$ws.Range("C4").Value = "Vicky"
$ws.Range("C3").Value = "Nico"

# Update the active selection to C4
$ws.Range("C4").Select() | Out-Null
